{"js": "// The edit updates the Java stack-trace text stored in the red/bold run of\n// the second paragraph: method/line numbers shift (the project moved from\n// Apache POI 3.16 to 3.17) and the JUnit stack frames gain a duplicated\n// Suite/ParentRunner call chain.\n//\n// The whole original stack-trace text is located with a single body search\n// and replaced in one operation, so the existing run/formatting (bold, red\n// font) on that run is left untouched.\n\nconst oldStackTrace =\n  \"divOp(java.lang.Integer,java.lang.Integer) with arguments [1, 0] failed:\\n\" +\n  \"\\t/ by zero\\n\" +\n  \"java.lang.ArithmeticException: / by zero\\n\" +\n  \"\\tat org.eclipse.acceleo.query.services.NumberServices.divOp(NumberServices.java:99)\\n\" +\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\" +\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\" +\n  \"\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\" +\n  \"\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\" +\n  \"\\tat org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:163)\\n\" +\n  \"\\tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:136)\\n\" +\n  \"\\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callService(EvaluationServices.java:129)\\n\" +\n  \"\\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:168)\\n\" +\n  \"\\tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:189)\\n\" +\n  \"\\tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:118)\\n\" +\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\" +\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\" +\n  \"\\tat org.eclipse.acceleo.query.parser.AstEvaluator.caseSequenceInExtensionLiteral(AstEvaluator.java:333)\\n\" +\n  \"\\tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:259)\\n\" +\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\" +\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\" +\n  \"\\tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112)\\n\" +\n  \"\\tat org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseRepetition(M2DocEvaluator.java:802)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseRepetition(M2DocEvaluator.java:1)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:139)\\n\" +\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\" +\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:836)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1034)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:183)\\n\" +\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\" +\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:836)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseTemplate(M2DocEvaluator.java:297)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseTemplate(M2DocEvaluator.java:1)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:201)\\n\" +\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\" +\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:836)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:259)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:246)\\n\" +\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\" +\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:836)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:252)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:691)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:396)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:318)\\n\" +\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\" +\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\" +\n  \"\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\" +\n  \"\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\" +\n  \"\\tat org.junit.runners.model.FrameworkMethod$1.runReflectiveCall(FrameworkMethod.java:50)\\n\" +\n  \"\\tat org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)\\n\" +\n  \"\\tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:47)\\n\" +\n  \"\\tat org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:325)\\n\" +\n  \"\\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:78)\\n\" +\n  \"\\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:57)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\\n\" +\n  \"\\tat org.junit.runners.Suite.runChild(Suite.java:128)\\n\" +\n  \"\\tat org.junit.runners.Suite.runChild(Suite.java:27)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\\n\" +\n  \"\\tat org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)\\n\" +\n  \"\\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\\n\" +\n  \"\\tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)\\n\" +\n  \"\\tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)\\n\" +\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:459)\\n\" +\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:675)\\n\" +\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:382)\\n\" +\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:192)\\n\" +\n  \"\";\n\nconst newStackTrace =\n  \"divOp(java.lang.Integer,java.lang.Integer) with arguments [1, 0] failed:\\n\" +\n  \"\\t/ by zero\\n\" +\n  \"java.lang.ArithmeticException: / by zero\\n\" +\n  \"\\tat org.eclipse.acceleo.query.services.NumberServices.divOp(NumberServices.java:99)\\n\" +\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\" +\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\" +\n  \"\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\" +\n  \"\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\" +\n  \"\\tat org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:163)\\n\" +\n  \"\\tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:136)\\n\" +\n  \"\\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callService(EvaluationServices.java:129)\\n\" +\n  \"\\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:168)\\n\" +\n  \"\\tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:189)\\n\" +\n  \"\\tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:118)\\n\" +\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\" +\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\" +\n  \"\\tat org.eclipse.acceleo.query.parser.AstEvaluator.caseSequenceInExtensionLiteral(AstEvaluator.java:333)\\n\" +\n  \"\\tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:259)\\n\" +\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\" +\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\" +\n  \"\\tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112)\\n\" +\n  \"\\tat org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseRepetition(M2DocEvaluator.java:1003)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseRepetition(M2DocEvaluator.java:1)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:147)\\n\" +\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\" +\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1254)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)\\n\" +\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\" +\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:275)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)\\n\" +\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\" +\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:264)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:712)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:459)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:369)\\n\" +\n  \"\\tat sun.reflect.GeneratedMethodAccessor75.invoke(Unknown Source)\\n\" +\n  \"\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\" +\n  \"\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\" +\n  \"\\tat org.junit.runners.model.FrameworkMethod$1.runReflectiveCall(FrameworkMethod.java:50)\\n\" +\n  \"\\tat org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)\\n\" +\n  \"\\tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:47)\\n\" +\n  \"\\tat org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)\\n\" +\n  \"\\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:325)\\n\" +\n  \"\\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:78)\\n\" +\n  \"\\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:57)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\\n\" +\n  \"\\tat org.junit.runners.Suite.runChild(Suite.java:128)\\n\" +\n  \"\\tat org.junit.runners.Suite.runChild(Suite.java:27)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\\n\" +\n  \"\\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\\n\" +\n  \"\\tat org.junit.runners.Suite.runChild(Suite.java:128)\\n\" +\n  \"\\tat org.junit.runners.Suite.runChild(Suite.java:27)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\\n\" +\n  \"\\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\\n\" +\n  \"\\tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)\\n\" +\n  \"\\tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)\\n\" +\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:539)\\n\" +\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:761)\\n\" +\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:461)\\n\" +\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:207)\\n\" +\n  \"\";\n\nconst results = context.document.body.search(oldStackTrace, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly 1 match for the stack trace text, found \" + results.items.length\n  );\n}\n\nresults.items[0].insertText(newStackTrace, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The edit updates the Java stack-trace text stored in the red/bold run of\n# the second paragraph: method/line numbers shift (the project moved from\n# Apache POI 3.16 to 3.17) and the JUnit stack frames gain a duplicated\n# Suite/ParentRunner call chain.\n#\n# The whole original stack-trace text is located with Word's Find and\n# replaced in a single pass, so the existing run formatting (bold, red font)\n# on that run is left untouched.\n\n$d = $word.ActiveDocument\n\n$oldStackTrace = @'\ndivOp(java.lang.Integer,java.lang.Integer) with arguments [1, 0] failed:\n\t/ by zero\njava.lang.ArithmeticException: / by zero\n\tat org.eclipse.acceleo.query.services.NumberServices.divOp(NumberServices.java:99)\n\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\n\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\n\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\n\tat java.lang.reflect.Method.invoke(Method.java:498)\n\tat org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:163)\n\tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:136)\n\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callService(EvaluationServices.java:129)\n\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:168)\n\tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:189)\n\tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:118)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\tat org.eclipse.acceleo.query.parser.AstEvaluator.caseSequenceInExtensionLiteral(AstEvaluator.java:333)\n\tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:259)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112)\n\tat org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseRepetition(M2DocEvaluator.java:802)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseRepetition(M2DocEvaluator.java:1)\n\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:139)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:836)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1034)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)\n\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:183)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:836)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseTemplate(M2DocEvaluator.java:297)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseTemplate(M2DocEvaluator.java:1)\n\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:201)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:836)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:259)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)\n\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:246)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:836)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:252)\n\tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:691)\n\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:396)\n\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:318)\n\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\n\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\n\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\n\tat java.lang.reflect.Method.invoke(Method.java:498)\n\tat org.junit.runners.model.FrameworkMethod$1.runReflectiveCall(FrameworkMethod.java:50)\n\tat org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)\n\tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:47)\n\tat org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)\n\tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:325)\n\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:78)\n\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:57)\n\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\n\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\n\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\n\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\n\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\n\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\n\tat org.junit.runners.Suite.runChild(Suite.java:128)\n\tat org.junit.runners.Suite.runChild(Suite.java:27)\n\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\n\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\n\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\n\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\n\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\n\tat org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)\n\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\n\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\n\tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)\n\tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)\n\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:459)\n\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:675)\n\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:382)\n\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:192)\n\n'@\n\n$newStackTrace = @'\ndivOp(java.lang.Integer,java.lang.Integer) with arguments [1, 0] failed:\n\t/ by zero\njava.lang.ArithmeticException: / by zero\n\tat org.eclipse.acceleo.query.services.NumberServices.divOp(NumberServices.java:99)\n\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\n\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\n\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\n\tat java.lang.reflect.Method.invoke(Method.java:498)\n\tat org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:163)\n\tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:136)\n\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callService(EvaluationServices.java:129)\n\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:168)\n\tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:189)\n\tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:118)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\tat org.eclipse.acceleo.query.parser.AstEvaluator.caseSequenceInExtensionLiteral(AstEvaluator.java:333)\n\tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:259)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112)\n\tat org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseRepetition(M2DocEvaluator.java:1003)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseRepetition(M2DocEvaluator.java:1)\n\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:147)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1254)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)\n\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:275)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)\n\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:264)\n\tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:712)\n\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:459)\n\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:369)\n\tat sun.reflect.GeneratedMethodAccessor75.invoke(Unknown Source)\n\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\n\tat java.lang.reflect.Method.invoke(Method.java:498)\n\tat org.junit.runners.model.FrameworkMethod$1.runReflectiveCall(FrameworkMethod.java:50)\n\tat org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)\n\tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:47)\n\tat org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)\n\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\n\tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:325)\n\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:78)\n\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:57)\n\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\n\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\n\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\n\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\n\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\n\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\n\tat org.junit.runners.Suite.runChild(Suite.java:128)\n\tat org.junit.runners.Suite.runChild(Suite.java:27)\n\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\n\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\n\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\n\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\n\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\n\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\n\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\n\tat org.junit.runners.Suite.runChild(Suite.java:128)\n\tat org.junit.runners.Suite.runChild(Suite.java:27)\n\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\n\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\n\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\n\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\n\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\n\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\n\tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)\n\tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)\n\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:539)\n\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:761)\n\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:461)\n\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:207)\n\n'@\n\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$result = $find.Execute($oldStackTrace, $false, $false, $false, $false, $false, $true, 1, $false, $newStackTrace, 2)\n\nif (-not $result) {\n    throw \"Could not find the expected stack trace text to replace.\"\n}\n"}
